# Remove the "Individual pages" report content (all ten paragraphs of
# page/file listings for Rob and Sinead), leaving only the single
# trailing empty paragraph that precedes the section properties.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
if ($count -gt 1) {
    # Everything from the start of the first paragraph through the end
    # of the second-to-last paragraph is the content to remove; the
    # final (already empty) paragraph is left untouched.
    $start = $d.Paragraphs.Item(1).Range.Start
    $end = $d.Paragraphs.Item($count - 1).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
